$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the date/time stamp in the document header.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "June  16, 2021 (08:22:28 PM)", $false, $false, $false, $false, $false,
    $true, 1, $false, "June  16, 2021 (08:24:03 PM)", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Append new sentences to the paragraph that ends with
#    "solution, that mixes classes and decision structures."
#
#    Two of the new runs ("what" / "how") must be italic. This engine
#    has no working API to add the <w:i/><w:iCs/> pair to brand-new
#    text via Font properties, but there happens to already be one
#    italic run in the document ("only then"). We borrow its
#    FormattedText (which really does carry the <w:i/><w:iCs/> run
#    properties), temporarily retarget its text to what we need,
#    paste a copy at our insertion point, then restore the original
#    run's text so the rest of the document is unaffected.
# ------------------------------------------------------------------

# Locate the italic template run ("only then") once, up front - its
# position never changes because every later edit happens further
# down in the document.
$tmplRng = $d.Content
$tmplRng.Find.Execute("only then") | Out-Null
$srcStart = $tmplRng.Start
$srcOriginalText = $tmplRng.Text
$srcOriginalLen = $srcOriginalText.Length

# Locate the insertion point (end of the target sentence).
$targetRng = $d.Content
$targetRng.Find.Execute("solution, that mixes classes and decision structures.") | Out-Null
$pos = $targetRng.End

# --- run: " " ---
$t = " "
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length

# --- run: "Spend some time reading through the implementation to understand" ---
$t = "Spend some time reading through the implementation to understand"
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length

# --- run: " " ---
$t = " "
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length

# --- run: "what"  (italic) ---
$newWord = "what"
$ft = $tmplRng.FormattedText
$ft.Text = $newWord
$shift = $newWord.Length - $srcOriginalLen
$pos = $pos + $shift
$d.Range($pos, $pos).FormattedText = $ft
$pos = $pos + $newWord.Length
$restoreRng = $d.Range($srcStart, $srcStart + $newWord.Length)
$restoreRng.Text = $srcOriginalText
$pos = $pos - $shift

# --- run: " " ---
$t = " "
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length

# --- run: "the program is doing and" ---
$t = "the program is doing and"
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length

# --- run: " " ---
$t = " "
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length

# --- run: "how"  (italic) ---
$newWord = "how"
$ft = $tmplRng.FormattedText
$ft.Text = $newWord
$shift = $newWord.Length - $srcOriginalLen
$pos = $pos + $shift
$d.Range($pos, $pos).FormattedText = $ft
$pos = $pos + $newWord.Length
$restoreRng = $d.Range($srcStart, $srcStart + $newWord.Length)
$restoreRng.Text = $srcOriginalText
$pos = $pos - $shift

# --- run: " " ---
$t = " "
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length

# --- run: "it is doing it." ---
$t = "it is doing it."
$d.Range($pos, $pos).InsertAfter($t) | Out-Null
$pos = $pos + $t.Length
